$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add the new "2022" column (column S) to the statistics table.
# ---------------------------------------------------------------------------

# 1. Copy the formatting that the neighbouring 2021 column (R) already uses
#    for every affected row into the new column S. This reuses the existing
#    number format (0.0), font, alignment (right / center) and, for the
#    last data row, the bottom border - instead of inventing new look and
#    feel for the new column.
$ws.Range("R4:R8").Copy() | Out-Null
$ws.Range("S4:S8").PasteSpecial(-4122) | Out-Null

$ws.Range("R10:R11").Copy() | Out-Null
$ws.Range("S10:S11").PasteSpecial(-4122) | Out-Null

$ws.Range("R13:R16").Copy() | Out-Null
$ws.Range("S13:S16").PasteSpecial(-4122) | Out-Null

$ws.Range("R18:R44").Copy() | Out-Null
$ws.Range("S18:S44").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# 2. Write the new "2022" header and data values into column S.
$ws.Range("S4").Value2 = 2022

$ws.Range("S5").Value2 = 33.152856050161155

$ws.Range("S7").Value2 = 32.831913512166025
$ws.Range("S8").Value2 = 33.509346380994529

$ws.Range("S10").Value2 = 34.041194942162896
$ws.Range("S11").Value2 = 32.636018013483323

$ws.Range("S13").Value2 = 40.271414365477746
$ws.Range("S14").Value2 = 31.568157010024336
$ws.Range("S15").Value2 = 30.277813022272248
$ws.Range("S16").Value2 = 22.733608300917229

$ws.Range("S18").Value2 = 48.492370829119814
$ws.Range("S19").Value2 = 46.987664282528065
$ws.Range("S20").Value2 = 50.118899291215271
$ws.Range("S21").Value2 = 47.142900749295329
$ws.Range("S22").Value2 = 47.199946558584017
$ws.Range("S23").Value2 = 47.082025761639336
$ws.Range("S24").Value2 = 31.240016364696597
$ws.Range("S25").Value2 = 31.228685777194666
$ws.Range("S26").Value2 = 31.252112297543153
$ws.Range("S27").Value2 = 42.049857693482664
$ws.Range("S28").Value2 = 42.689244289315013
$ws.Range("S29").Value2 = 41.380596558931735
$ws.Range("S30").Value2 = 19.945481087558658
$ws.Range("S31").Value2 = 20.799187962023481
$ws.Range("S32").Value2 = 19.013188474520234
$ws.Range("S33").Value2 = 23.919779113642239
$ws.Range("S34").Value2 = 23.962040711070269
$ws.Range("S35").Value2 = 23.876854008981983
$ws.Range("S36").Value2 = 26.113584517813127
$ws.Range("S37").Value2 = 25.651528441631889
$ws.Range("S38").Value2 = 26.620973515499056
$ws.Range("S39").Value2 = 35.676666099583812
$ws.Range("S40").Value2 = 34.026766685280904
$ws.Range("S41").Value2 = 37.792274390474752
$ws.Range("S42").Value2 = 26.602385500795538
$ws.Range("S43").Value2 = 25.585637135242425
$ws.Range("S44").Value2 = 27.750206810614948

# Row 6 only needs the (blank) formatting applied above - no value to write.

# 3. The two subtotal rows (9 and 12) must not gain a 2022 column at all -
#    clear both the (inherited) formatting and any content from those cells.
$ws.Range("S9").ClearContents() | Out-Null
$ws.Range("S9").ClearFormats() | Out-Null
$ws.Range("S12").ClearContents() | Out-Null
$ws.Range("S12").ClearFormats() | Out-Null

# ---------------------------------------------------------------------------
# Update the active selection to reflect where the editor left off.
# ---------------------------------------------------------------------------
$ws.Range("R8").Select() | Out-Null
